# TC02_Canine_Filter_FileType-IndexFile.xlsx
# The "CasesTab" Cypher query (cell B2 on the "startup" sheet) included a
# trailing OPTIONAL MATCH / coalesce(...) clause that pulled back a `Cohort`
# column which isn't actually surfaced by the report. Drop that clause so the
# query only returns the columns that are really used (Diagnosis, FileType,
# FileFormat, FileAssociation, NeuteredStatus, PrimeDiseaseSite, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
  "MATCH (c)<--(diag:diagnosis)`n" +
  "MATCH (samp:sample)-->(c) `n" +
  "  MATCH (f:file)-[*]->(c)`n" +
  "   WHERE f.file_type IN [""Index File""] `n" +
  "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
  "  WITH DISTINCT c, s, demo, diag, co`n" +
  "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
  "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
  "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
  "        coalesce(demo.breed, '') AS Breed ,`n" +
  "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
  "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
  "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
  "        coalesce(demo.sex, '') AS Sex ,`n" +
  "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
  "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
  "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# The shorter query text re-wraps onto fewer lines, so the row shrinks a bit.
$ws.Rows(2).RowHeight = 259.2
$ws.Rows(3).RowHeight = 288
$ws.Rows(4).RowHeight = 259.2

# Leave the selection on the cell that was actually edited.
$ws.Range("B2").Select()
